$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.770.86"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").Value = "1.748.82"
$ws.Range("E3").Value = "  -5.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5023"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "41.72"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2637"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -13.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06146"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.78%  "
$ws.Range("D11").Value = "1.750.74"
$ws.Range("E11").Value = "  -5.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06945"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -15.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.485"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -10.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5883"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -20.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.84"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -14.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "25.819.52"
$ws.Range("E19").Value = "  -2.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.66"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -16.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006761"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -14.91%  "
$ws.Range("D22").Value = "1.971.50"
$ws.Range("E22").Value = "  -5.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.058"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.102"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -13.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.116"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.535"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.30%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.827"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -17.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "14.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.20%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.754"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -11.93%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08109"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.460"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -14.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04486"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.639"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9839"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -13.30%  "
$ws.Range("E38").Value = "  -16.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.662"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01542"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.31%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.916"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -16.77%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "104.32"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.146"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -12.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3797"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -20.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7244"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -20.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05304"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.52%  "
$ws.Range("E48").Value = "  -11.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -14.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.904"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -21.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -13.61%  "
